# Update loading_percent values for the 380 kV case (Case_4_10).
# Columns B,C,D,E,G,H,I,J,L,N are rewritten for data rows 2-25;
# columns F,K,M,O stay at 0 and column A (index) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B", "C", "D", "E", "G", "H", "I", "J", "L", "N")

# Each inner array holds the new values for columns B,C,D,E,G,H,I,J,L,N,
# in that order, for sheet rows 2 through 25.
$newValues = @(
    @(24.23365319659188,9.43623903500807,12.87209423634133,12.5273083942387,57.96440946188886,21.41963293469,32.34717319489383,7.759257348099688,13.87585490966716,20.11759498332341),
    @(23.89247543999452,9.034822717385589,12.88138324503198,12.55492113915233,57.85729736146605,21.4561897580185,32.42360694753853,7.758821537352147,13.87212774838779,20.19232851165093),
    @(23.68629876282057,8.781396656597179,12.88933671721171,12.57313491652429,57.80877182782097,21.4830778761529,32.47760160799495,7.758608764314492,13.87222642756381,20.24020280808434),
    @(23.60321467199493,8.676541264259358,12.89314348222514,12.58087436286751,57.79333603968338,21.49514907700746,32.5013750873497,7.758535841523546,13.87286787631525,20.2602136138299),
    @(23.58947804254973,8.659039811088157,12.89380975058731,12.58217866356745,57.79103499155867,21.49722068777143,32.50542938869462,7.758524565271019,13.8730107245549,20.26356674479214),
    @(23.68517434301214,8.779988711359342,12.88938576650964,12.57323800848752,57.80854608735096,21.4832361658786,32.47791506628877,7.758607725039661,13.87223264291655,20.2404706472893),
    @(24.11539063167996,9.299367966786628,12.87483007720729,12.53656816474569,57.92390197579967,21.43131408511101,32.37205783109874,7.759095702843998,13.8740749779393,20.14295169043609),
    @(24.98072505451243,10.25635514517077,12.86413922800951,12.47463011593616,58.28657524292148,21.36486311766885,32.22078449610242,7.760487775990218,13.89657753144136,19.96740247423269),
    @(25.62358602438325,10.91475053106018,12.8671611239164,12.43517351206187,58.63546304878776,21.33775714504614,32.14431762631899,7.761776387787359,13.92454350081534,19.84786641366795),
    @(25.9163033614521,11.20336999273585,12.87089257801021,12.41853103853058,58.8118212548215,21.33016872009316,32.1171271558448,7.762420485991645,13.93972550511064,19.79550992275061),
    @(26.02708360260816,11.3110162684026,12.87264370058922,12.41241637544836,58.88111181519331,21.32797880663199,32.1079278334711,7.762672722467371,13.94582582950173,19.77597254267955),
    @(26.00322965054428,11.28790717202763,12.87225154260409,12.41372494558297,58.86607785722846,21.32842001354868,32.10986020178239,7.762618028260106,13.94449643787846,19.78016744634205),
    @(25.92541922265396,11.21225956395454,12.87102987311166,12.41802422639034,58.81747176469646,21.32997484426881,32.11634830683986,7.762441070657903,13.94022035884974,19.79389679095485),
    @(25.87774649122048,11.16570642429924,12.8703255691297,12.42068206380045,58.78802471381373,21.33101630149442,32.12046547789706,7.762333764019592,13.93764679725113,19.80234398102623),
    @(25.60445276208107,10.89566218421843,12.86696462498597,12.43628739355279,58.62429026231212,21.33834864738351,32.14624784924665,7.761735459947705,13.92360061509464,19.851328553553),
    @(25.43679098628152,10.72714778444065,12.86550591327937,12.44619510809074,58.52834439065246,21.34406273297418,32.1640133682178,7.761383252822846,13.91561212381071,19.88189540466576),
    @(25.34038853955555,10.62920198981162,12.86488879199569,12.45201676404613,58.47482201459445,21.34779560016155,32.174946197564,7.76118612236514,13.91124908519125,19.89966698453852),
    @(25.30775710693203,10.59586661292076,12.86471797658683,12.45400901437385,58.45698674984826,21.34913607518537,32.1787704352262,7.761120313259245,13.90981170784202,19.90571687953023),
    @(25.4546363106494,10.7451927089479,12.86563823595396,12.44512768778163,58.53838606710583,21.34340825753504,32.16204820558388,7.761420181810502,13.9164385466614,19.87862182490372),
    @(25.94827663191816,11.23452440945133,12.87137953838649,12.41675633837927,58.83168076083252,21.32949958626067,32.11441277959464,7.762492821274226,13.94146683627403,19.7898563231596),
    @(26.27047097139868,11.54469753567673,12.87710202237926,12.39930662002556,59.03796453796431,21.32439492821003,32.08967717166887,7.763242411656261,13.95987012668572,19.73352598731074),
    @(26.098582932229,11.38005773619246,12.87386786402476,12.40852002031763,58.92654231986856,21.32675421985647,32.10229215935333,7.762837896593208,13.94986165871986,19.76343712128445),
    @(25.44656846711068,10.73703790833321,12.8655777228478,12.44560987731505,58.53384112014835,21.34370275123884,32.16293441572552,7.761403469527539,13.91606420537131,19.88010119249982),
    @(24.74494235413925,10.00485978997667,12.86512014231033,12.49032156565031,58.17391302280183,21.37903643645789,32.25564249555946,7.760064538808736,13.88847511389357,20.01322663798352)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $newValues[$i]
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $ws.Range("$($columns[$j])$rowNum").Value = $rowValues[$j]
    }
}
